$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'69.806.02"
$ws.Cells.Item(2, 5).Value = '  -1.45%  '
$ws.Cells.Item(3, 4).Value = "'3.514.23"
$ws.Cells.Item(3, 5).Value = '  -1.48%  '
$ws.Cells.Item(4, 5).Value = '  -0.22%  '
$ws.Cells.Item(5, 4).Value = "'616.02"
$ws.Cells.Item(6, 4).Value = "'192.15"
$ws.Cells.Item(6, 5).Value = '  +1.20%  '
$ws.Cells.Item(7, 5).Value = '  +0.77%  '
$ws.Cells.Item(8, 5).Value = '  -0.11%  '
$ws.Cells.Item(9, 4).Value = "'0.213"
$ws.Cells.Item(9, 5).Value = '  -3.27%  '
$ws.Cells.Item(10, 4).Value = "'0.663"
$ws.Cells.Item(10, 5).Value = '  +1.92%  '
$ws.Cells.Item(11, 4).Value = "'53.41"
$ws.Cells.Item(11, 5).Value = '  -2.36%  '
$ws.Cells.Item(12, 5).Value = '  -3.68%  '
$ws.Cells.Item(13, 4).Value = "'9.58"
$ws.Cells.Item(13, 5).Value = '  +0.48%  '
$ws.Cells.Item(14, 4).Value = "'4.082.12"
$ws.Cells.Item(14, 5).Value = '  -1.35%  '
$ws.Cells.Item(15, 4).Value = "'621.27"
$ws.Cells.Item(15, 5).Value = '  +9.59%  '
$ws.Cells.Item(16, 4).Value = "'69.887.07"
$ws.Cells.Item(16, 5).Value = '  -1.38%  '
$ws.Cells.Item(17, 4).Value = "'19.01"
$ws.Cells.Item(17, 5).Value = '  -1.01%  '
$ws.Cells.Item(18, 5).Value = '  -1.09%  '
$ws.Cells.Item(19, 4).Value = "'3.517.04"
$ws.Cells.Item(19, 5).Value = '  -0.48%  '
$ws.Cells.Item(20, 5).Value = '  -0.29%  '
$ws.Cells.Item(21, 5).Value = '  -1.27%  '
$ws.Cells.Item(22, 4).Value = "'109.14"
$ws.Cells.Item(22, 5).Value = '  +15.75%  '
$ws.Cells.Item(23, 4).Value = "'17.18"
$ws.Cells.Item(23, 5).Value = '  -4.27%  '
$ws.Cells.Item(24, 5).Value = '  +2.11%  '
$ws.Cells.Item(25, 4).Value = "'5.03"
$ws.Cells.Item(25, 5).Value = '  +2.37%  '
$ws.Cells.Item(26, 4).Value = "'3.12"
$ws.Cells.Item(26, 5).Value = '  +6.37%  '
$ws.Cells.Item(27, 4).Value = "'10.96"
$ws.Cells.Item(27, 5).Value = '  -2.16%  '
$ws.Cells.Item(28, 4).Value = "'9.69"
$ws.Cells.Item(28, 5).Value = '  +4.19%  '
$ws.Cells.Item(29, 4).Value = "'34.31"
$ws.Cells.Item(29, 5).Value = '  +5.46%  '
$ws.Cells.Item(30, 4).Value = "'6.97"
$ws.Cells.Item(30, 5).Value = '  -3.53%  '
$ws.Cells.Item(31, 4).Value = "'12.54"
$ws.Cells.Item(31, 5).Value = '  +1.57%  '
$ws.Cells.Item(32, 2).Value = 'dogwifhat'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(32, 4).Value = "'3.94"
$ws.Cells.Item(32, 5).Value = '  +5.27%  '
$ws.Cells.Item(33, 2).Value = 'Hedera'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(33, 4).Value = "'0.116"
$ws.Cells.Item(33, 5).Value = '  -0.02%  '
$ws.Cells.Item(34, 4).Value = "'63.45"
$ws.Cells.Item(34, 5).Value = '  -0.98%  '
$ws.Cells.Item(35, 4).Value = "'3.10"
$ws.Cells.Item(35, 5).Value = '  -5.15%  '
$ws.Cells.Item(36, 2).Value = 'Bittensor'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(36, 4).Value = "'522.21"
$ws.Cells.Item(36, 5).Value = '  -2.05%  '
$ws.Cells.Item(37, 4).Value = "'3.664.97"
$ws.Cells.Item(37, 5).Value = '  +0.29%  '
$ws.Cells.Item(38, 2).Value = 'Dai'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(38, 4).Value = "'0.998"
$ws.Cells.Item(38, 5).Value = '  -0.08%  '
$ws.Cells.Item(39, 4).Value = "'3.64"
$ws.Cells.Item(39, 5).Value = '  +5.93%  '
$ws.Cells.Item(40, 5).Value = '  -4.44%  '
$ws.Cells.Item(41, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(41, 4).Value = "'36.68"
$ws.Cells.Item(41, 5).Value = '  -4.49%  '
$ws.Cells.Item(42, 2).Value = 'PEPE'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(42, 4).Value = "'0.0₃0779"
$ws.Cells.Item(42, 5).Value = '  -2.39%  '
$ws.Cells.Item(43, 5).Value = '  -1.41%  '
$ws.Cells.Item(44, 4).Value = "'0.0470"
$ws.Cells.Item(44, 5).Value = '  +0.06%  '
$ws.Cells.Item(45, 5).Value = '  -0.53%  '
$ws.Cells.Item(46, 4).Value = "'0.142"
$ws.Cells.Item(46, 5).Value = '  +3.04%  '
$ws.Cells.Item(47, 4).Value = "'3.32"
$ws.Cells.Item(47, 5).Value = '  -3.87%  '
$ws.Cells.Item(48, 5).Value = '  -5.53%  '
$ws.Cells.Item(49, 5).Value = '  +0.40%  '
$ws.Cells.Item(50, 4).Value = "'132.38"
$ws.Cells.Item(50, 5).Value = '  -1.19%  '
$ws.Cells.Item(51, 4).Value = "'0.000240"
$ws.Cells.Item(51, 5).Value = '  -4.77%  '

Write-Host "Updated 98 cells"